$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "f"
$ws.Range("B3").Value = "df"

$ws.Range("D5").Select()
